# Apply the edits described by the commit diff:
#  - Update column E ("time" / quantity?) values for rows 79-93 on Sheet1.
#    Rows 79-85 and 93 go from 8 -> 22; rows 86-92 go from 8 -> 23.
#  - Update the current selection on Sheet1 from G71 to J84.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update cell values in column E ---
$rows22 = @(79, 80, 81, 82, 83, 84, 85, 93)
foreach ($r in $rows22) {
    $ws.Cells.Item($r, 5).Value = 22
}

$rows23 = @(86, 87, 88, 89, 90, 91, 92)
foreach ($r in $rows23) {
    $ws.Cells.Item($r, 5).Value = 23
}

# --- Update the active selection to J84 (was G71) ---
$ws.Activate()
$ws.Range("J84").Select()
